$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.806.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'3.541.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.58%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'615.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -4.62%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'154.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.15%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.536.73"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -3.55%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.92%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.32%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.57%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.65%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -2.29%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'32.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.10%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.142.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.56%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.529.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.38%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'67.766.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.47%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +0.30%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.47%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'15.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.83%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'454.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.37%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'9.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.04%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.00%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'77.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.77%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.684.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.07%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -4.26%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.72%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'8.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.70%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.59%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.50%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.05%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -2.33%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -4.04%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -2.98%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D37").Value = "'3.541.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.34%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -3.88%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.06%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'176.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.45%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -4.44%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.94%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.35%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.888"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.46%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'29.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +10.18%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'45.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.77%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -3.27%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.28%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.80%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -3.55%  "
$ws.Range("E51").Style = "Normal"
